$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format so numeric-looking strings (with trailing zeros,
# multiple dots, etc.) are preserved exactly as typed instead of being reinterpreted
# as numbers.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "28.661.56"
$ws.Range("E2").Value = "  -1.82%  "

# Row 3
$ws.Range("D3").Value = "1.802.11"
$ws.Range("E3").Value = "  -1.30%  "

# Row 4
$ws.Range("E4").Value = "  +0.19%  "

# Row 5
$ws.Range("D5").Value = "232.00"
$ws.Range("E5").Value = "  -0.89%  "

# Row 6
$ws.Range("D6").Value = "0.5921"
$ws.Range("E6").Value = "  -1.55%  "

# Row 7
$ws.Range("E7").Value = "  +0.19%  "

# Row 8
$ws.Range("D8").Value = "0.2774"
$ws.Range("E8").Value = "  -0.63%  "

# Row 9
$ws.Range("D9").Value = "0.06815"
$ws.Range("E9").Value = "  -3.39%  "

# Row 10
$ws.Range("D10").Value = "23.30"

# Row 11
$ws.Range("B11").Value = "TRON"
$ws.Range("C11").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D11").Value = "0.07514"
$ws.Range("E11").Value = "  -1.60%  "

# Row 12
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.799.88"
$ws.Range("E12").Value = "  -1.38%  "

# Row 13
$ws.Range("D13").Value = "4.765"
$ws.Range("E13").Value = "  -0.27%  "

# Row 14
$ws.Range("D14").Value = "0.6223"
$ws.Range("E14").Value = "  -0.51%  "

# Row 15
$ws.Range("D15").Value = "2.046.44"
$ws.Range("E15").Value = "  -1.31%  "

# Row 16
$ws.Range("D16").Value = "0.000009170"
$ws.Range("E16").Value = "  -7.60%  "

# Row 17
$ws.Range("E17").Value = "  -4.16%  "

# Row 18
$ws.Range("D18").Value = "28.644.67"
$ws.Range("E18").Value = "  -1.82%  "

# Row 19
$ws.Range("D19").Value = "5.462"
$ws.Range("E19").Value = "  -6.35%  "

# Row 20
$ws.Range("E20").Value = "  +0.20%  "

# Row 21
$ws.Range("D21").Value = "210.64"
$ws.Range("E21").Value = "  -6.79%  "

# Row 22
$ws.Range("D22").Value = "11.50"
$ws.Range("E22").Value = "  -1.50%  "

# Row 23
$ws.Range("D23").Value = "6.831"
$ws.Range("E23").Value = "  -2.24%  "

# Row 24
$ws.Range("E24").Value = "  +0.27%  "

# Row 25
$ws.Range("D25").Value = "153.80"
$ws.Range("E25").Value = "  -1.11%  "

# Row 26
$ws.Range("D26").Value = "7.850"
$ws.Range("E26").Value = "  -1.98%  "

# Row 27
$ws.Range("D27").Value = "0.1267"
$ws.Range("E27").Value = "  -2.28%  "

# Row 28
$ws.Range("D28").Value = "16.45"
$ws.Range("E28").Value = "  -0.58%  "

# Row 29
$ws.Range("D29").Value = "1.405"
$ws.Range("E29").Value = "  -4.56%  "

# Row 30
$ws.Range("D30").Value = "0.06159"
$ws.Range("E30").Value = "  -0.63%  "

# Row 31
$ws.Range("D31").Value = "1.424"
$ws.Range("E31").Value = "  -1.34%  "

# Row 32
$ws.Range("D32").Value = "3.773"
$ws.Range("E32").Value = "  -1.28%  "

# Row 33
$ws.Range("D33").Value = "3.740"
$ws.Range("E33").Value = "  -1.26%  "

# Row 34
$ws.Range("D34").Value = "1.727"
$ws.Range("E34").Value = "  -0.85%  "

# Row 35
$ws.Range("D35").Value = "1.057"
$ws.Range("E35").Value = "  -5.51%  "

# Row 36
$ws.Range("D36").Value = "0.6417"
$ws.Range("E36").Value = "  +0.51%  "

# Row 37
$ws.Range("D37").Value = "2.501"
$ws.Range("E37").Value = "  -1.35%  "

# Row 38
$ws.Range("E38").Value = "  -0.54%  "

# Row 39
$ws.Range("D39").Value = "6.548"
$ws.Range("E39").Value = "  +0.31%  "

# Row 40
$ws.Range("D40").Value = "0.01690"
$ws.Range("E40").Value = "  -2.98%  "

# Row 41
$ws.Range("D41").Value = "1.149.74"
$ws.Range("E41").Value = "  -5.46%  "

# Row 42
$ws.Range("D42").Value = "0.8840"
$ws.Range("E42").Value = "  -1.77%  "

# Row 43
$ws.Range("D43").Value = "1.007"
$ws.Range("E43").Value = "  +0.36%  "

# Row 44
$ws.Range("D44").Value = "100.15"
$ws.Range("E44").Value = "  -0.18%  "

# Row 45
$ws.Range("D45").Value = "1.951.87"
$ws.Range("E45").Value = "  -1.64%  "

# Row 46
$ws.Range("D46").Value = "60.41"
$ws.Range("E46").Value = "  -3.42%  "

# Row 47
$ws.Range("D47").Value = "0.00000000113"
$ws.Range("E47").Value = "  -1.95%  "

# Row 48
$ws.Range("D48").Value = "1.585"
$ws.Range("E48").Value = "  +0.16%  "

# Row 49
$ws.Range("D49").Value = "8.370"
$ws.Range("E49").Value = "  -1.74%  "

# Row 50
$ws.Range("E50").Value = "  -0.71%  "

# Row 51
$ws.Range("E51").Value = "  -1.76%  "
